# --- edit.ps1 ---
# Updates the "map" sheet of Procedure.xlsx to the 8-7-2025 XtEHR version:
# header/presentedForm/bodySite/reason[x] restructuring of the EHDSProcedure rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old row layout (rows 2-23), keep the bold header-row style on row 1 intact.
$ws.Range("A2:D23").ClearContents()

# --- Row 1 (column headers) ---
$ws.Range("A1").Value = 'xtehr'
$ws.Range("B1").Value = 'zib'

# --- Rebuild rows 2-32 with the new mapping content ---
# row 2
$ws.Range("A2").Value = 'EHDSProcedure'
$ws.Range("B2").Value = 'Procedure'
$ws.Range("C2").Value = 'x'
$ws.Range("D2").Value = 'x'

# row 3
$ws.Range("A3").Value = 'EHDSProcedure.header'
$ws.Range("B3").Value = 'Registratiegegevens, zie DataSet mapping'

# row 4
$ws.Range("A4").Value = 'EHDSProcedure.header.subject'

# row 5
$ws.Range("A5").Value = 'EHDSProcedure.header.identifier'

# row 6
$ws.Range("A6").Value = 'EHDSProcedure.header.authorship'

# row 7
$ws.Range("A7").Value = 'EHDSProcedure.header.authorship.author[x]'

# row 8
$ws.Range("A8").Value = 'EHDSProcedure.header.authorship.datetime'

# row 9
$ws.Range("A9").Value = 'EHDSProcedure.header.lastUpdate'

# row 10
$ws.Range("A10").Value = 'EHDSProcedure.header.status'

# row 11
$ws.Range("A11").Value = 'EHDSProcedure.header.statusReason[x]'

# row 12
$ws.Range("A12").Value = 'EHDSProcedure.header.language'

# row 13
$ws.Range("A13").Value = 'EHDSProcedure.header.version'

# row 14
$ws.Range("A14").Value = 'EHDSProcedure.presentedForm'

# row 15
$ws.Range("A15").Value = 'EHDSProcedure.code'
$ws.Range("B15").Value = 'Procedure.ProcedureType'
$ws.Range("C15").Value = 'x'
$ws.Range("D15").Value = 'x'

# row 16
$ws.Range("A16").Value = 'EHDSProcedure.date[x]'
$ws.Range("B16").Value = 'Procedure.ProcedureStartDate'
$ws.Range("C16").Value = 'x'
$ws.Range("D16").Value = 'x'

# row 17
$ws.Range("A17").Value = 'EHDSProcedure.date[x]'
$ws.Range("B17").Value = 'Procedure.ProcedureEndDate'
$ws.Range("C17").Value = 'x'
$ws.Range("D17").Value = 'x'

# row 18
$ws.Range("A18").Value = 'EHDSProcedure.performer'
$ws.Range("B18").Value = 'Procedure.Performer::HealthProfessional'
$ws.Range("C18").Value = 'x'
$ws.Range("D18").Value = 'x'

# row 19
$ws.Range("B19").Value = 'Procedure.Requester::HealthProfessional'
$ws.Range("C19").Value = 'x'

# row 20
$ws.Range("A20").Value = 'EHDSProcedure.bodySite'
$ws.Range("B20").Value = 'Procedure.ProcedureAnatomicalLocation::AnatomicalLocation'
$ws.Range("C20").Value = 'x'
$ws.Range("D20").Value = 'x'

# row 21
$ws.Range("A21").Value = 'EHDSProcedure.reason[x]'
$ws.Range("B21").Value = 'Procedure.Indication::Problem'
$ws.Range("C21").Value = 'x'

# row 22
$ws.Range("A22").Value = 'EHDSProcedure.outcome'

# row 23
$ws.Range("A23").Value = 'EHDSProcedure.complication'

# row 24
$ws.Range("A24").Value = 'EHDSProcedure.deviceUsed'

# row 25
$ws.Range("B25").Value = 'Procedure.ProcedureMethod'
$ws.Range("C25").Value = 'x'
$ws.Range("D25").Value = 'x'

# row 26
$ws.Range("A26").Value = 'EHDSProcedure.focalDevice'
$ws.Range("B26").Value = 'Procedure.MedicalDevice'
$ws.Range("C26").Value = 'x'
$ws.Range("D26").Value = 'x'

# row 27
$ws.Range("A27").Value = 'EHDSProcedure.location'
$ws.Range("B27").Value = 'Procedure.Location::HealthcareProvider'
$ws.Range("C27").Value = 'x'
$ws.Range("D27").Value = 'x'

# row 28
$ws.Range("A28").Value = 'EHDSProcedure.note'

# row 29
$ws.Range("A29").Value = 'EHDSProcedure.reason[x]'
$ws.Range("B29").Value = 'Procedure.Indication'
$ws.Range("D29").Value = 'x'

# row 30
$ws.Range("A30").Value = 'EHDSProcedure.reason[x]'
$ws.Range("B30").Value = 'Procedure.Indication.Diagnosis'
$ws.Range("D30").Value = 'x'

# row 31
$ws.Range("A31").Value = 'EHDSProcedure.reason[x]'
$ws.Range("B31").Value = 'Procedure.Indication.Reaction'
$ws.Range("D31").Value = 'x'

# row 32
$ws.Range("A32").Value = 'EHDSProcedure.reason[x]'
$ws.Range("B32").Value = 'Procedure.Indication.Symptom'
$ws.Range("D32").Value = 'x'

# --- Column A widened (longer "EHDSProcedure.header..." labels) ---
# Target best-fit width from the authored workbook is 41.5703125 (computed by
# real Excel from actual glyph metrics). This runtime's ColumnWidth setter
# quantizes to coarser steps, so 40.6 is the closest request that lands on
# the nearest representable width (41.5).
$ws.Columns("A").ColumnWidth = 40.6

# --- Selection moved to A14 (EHDSProcedure.presentedForm) ---
$ws.Range("A14").Select()
